# "TAKEN 1 & 2" — add Marking.Notes / Notes.Format / Feedback.to.Learner /
# Feedback.Format columns (from a re-exported grade-centre "user info with
# coding" report) and refresh the embedded score-column header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The score header encodes a running submission id; it moved from
# ...127326 to ...130753 in the re-export.
$ws.Range("G1").Value = "TASK0.instructions..Total.Pts..1.Score...130753"

# Insert 4 new columns right after the score column (H:K), pushing the
# existing "Group Code..newid" block from H:M to L:Q.
$ws.Range("H1:K1").EntireColumn.Insert()

# New column headers; the data rows underneath are left blank (no marking
# notes/feedback recorded yet).
$ws.Range("H1").Value = "Marking.Notes"
$ws.Range("I1").Value = "Notes.Format"
$ws.Range("J1").Value = "Feedback.to.Learner"
$ws.Range("K1").Value = "Feedback.Format"

# Row 4 (c1243957 / Leonard Maaya) previously had no Student.ID on file;
# the re-export fills it in.
$ws.Range("D4").Value = 31243957
